$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "302.25"
Set-TextValue $ws.Range("E2") "2.32%"

Set-TextValue $ws.Range("D3") "42.54"
Set-TextValue $ws.Range("E3") "5.01%"

Set-TextValue $ws.Range("D4") "5.027"
Set-TextValue $ws.Range("E4") "0.21%"

Set-TextValue $ws.Range("D5") "0.07679"
Set-TextValue $ws.Range("E5") "3.83%"

Set-TextValue $ws.Range("D6") "1.609"
Set-TextValue $ws.Range("E6") "2.41%"

Set-TextValue $ws.Range("D7") "0.9970"
Set-TextValue $ws.Range("E7") "7.82%"

Set-TextValue $ws.Range("E8") "0.31%"

Set-TextValue $ws.Range("D9") "0.1208"
Set-TextValue $ws.Range("E9") "-0.77%"

Set-TextValue $ws.Range("D10") "0.1850"
Set-TextValue $ws.Range("E10") "2.03%"

Set-TextValue $ws.Range("D11") "0.09081"
Set-TextValue $ws.Range("E11") "3.34%"

Set-TextValue $ws.Range("D12") "0.04104"
Set-TextValue $ws.Range("E12") "-6.70%"

Set-TextValue $ws.Range("D13") "0.1045"
Set-TextValue $ws.Range("E13") "-0.97%"

Set-TextValue $ws.Range("D14") "0.001266"
Set-TextValue $ws.Range("E14") "0.09%"

Set-TextValue $ws.Range("D15") "0.005984"
Set-TextValue $ws.Range("E15") "1.05%"

Set-TextValue $ws.Range("D16") "0.007430"
Set-TextValue $ws.Range("E16") "1,896.26%"

Set-TextValue $ws.Range("D17") "3.315"
Set-TextValue $ws.Range("E17") "-0.86%"

Set-TextValue $ws.Range("D18") "4.392"
Set-TextValue $ws.Range("E18") "2.32%"

Set-TextValue $ws.Range("D19") "0.3340"
Set-TextValue $ws.Range("E19") "-0.22%"

Set-TextValue $ws.Range("D20") "8.362"
Set-TextValue $ws.Range("E20") "5.99%"

Set-TextValue $ws.Range("D21") "0.1361"
Set-TextValue $ws.Range("E21") "-2.08%"

Set-TextValue $ws.Range("D22") "0.2985"
Set-TextValue $ws.Range("E22") "6.48%"

Set-TextValue $ws.Range("D23") "0.04143"
Set-TextValue $ws.Range("E23") "5.29%"

Set-TextValue $ws.Range("D24") "0.001265"
Set-TextValue $ws.Range("E24") "0.46%"

Set-TextValue $ws.Range("D25") "0.003952"
Set-TextValue $ws.Range("E25") "4.00%"

Set-TextValue $ws.Range("D26") "0.0001344"
Set-TextValue $ws.Range("E26") "9.31%"

Set-TextValue $ws.Range("D38") "0.02458"
Set-TextValue $ws.Range("E38") "5.29%"

Set-TextValue $ws.Range("D39") "0.05275"
Set-TextValue $ws.Range("E39") "3.71%"

Set-TextValue $ws.Range("D40") "0.005773"
Set-TextValue $ws.Range("E40") "-6.40%"

Set-TextValue $ws.Range("D41") "0.007637"
Set-TextValue $ws.Range("E41") "-2.43%"

Set-TextValue $ws.Range("D42") "0.1349"
Set-TextValue $ws.Range("E42") "4.24%"

Set-TextValue $ws.Range("D43") "0.007362"
Set-TextValue $ws.Range("E43") "0.08%"

Set-TextValue $ws.Range("D44") "0.007333"
Set-TextValue $ws.Range("E44") "-0.15%"

Set-TextValue $ws.Range("D45") "0.3028"
Set-TextValue $ws.Range("E45") "3.24%"

Set-TextValue $ws.Range("D46") "0.00006605"
Set-TextValue $ws.Range("E46") "5.93%"

Set-TextValue $ws.Range("D47") "0.00000000746"
Set-TextValue $ws.Range("E47") "-0.48%"

Set-TextValue $ws.Range("D48") "0.04594"
Set-TextValue $ws.Range("E48") "-1.99%"

Set-TextValue $ws.Range("E49") "0.01%"

Set-TextValue $ws.Range("D50") "0.00002090"
Set-TextValue $ws.Range("E50") "-0.48%"

Set-TextValue $ws.Range("D51") "0.0001991"
Set-TextValue $ws.Range("E51") "-0.48%"
